$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Column-level base formatting (mirrors the <cols> defaults in the target) ----
$ws.Columns("A").ColumnWidth = 13.57642857142857
$ws.Columns("B").ColumnWidth = 13.57642857142857
$ws.Columns("C").ColumnWidth = 13.57642857142857
$ws.Columns("D").ColumnWidth = 13.57642857142857

$ws.Columns("B").NumberFormat = "#,##0"
$ws.Columns("C").NumberFormat = "#,##0"
$ws.Columns("C").HorizontalAlignment = -4152

# ---- Row heights ----
$ws.Rows(1).RowHeight = 18.75
$ws.Rows(2).RowHeight = 19.5
$ws.Rows(3).RowHeight = 19.5
$ws.Rows(4).RowHeight = 19.5
$ws.Rows(5).RowHeight = 19.5
$ws.Rows(6).RowHeight = 18.75

# ---- Row 1: Lampada do quarto ----
$ws.Range("A1").Value = "Lampada do quarto"
$ws.Range("B1").Value = 100
$ws.Range("B1").NumberFormat = "#,##0"
$ws.Range("B1").HorizontalAlignment = -4152
$ws.Range("B1").Font.Color = 0
$ws.Range("C1").Value = $true
$ws.Range("C1").NumberFormat = "#,##0"
$ws.Range("C1").HorizontalAlignment = -4152
$ws.Range("C1").Font.ThemeColor = 1
$ws.Range("D1").NumberFormat = "General"

# ---- Row 2: Ar da sala ----
$ws.Range("A2").Value = "Ar da sala"
$ws.Range("B2").Value = 18
$ws.Range("B2").NumberFormat = "#,##0"
$ws.Range("B2").HorizontalAlignment = -4152
$ws.Range("B2").Font.Color = 0
$ws.Range("C2").Value = $true
$ws.Range("C2").NumberFormat = "#,##0"
$ws.Range("C2").HorizontalAlignment = -4152
$ws.Range("C2").Font.ThemeColor = 1
$ws.Range("D2").NumberFormat = "General"

# ---- Row 3: Tv do quarto ----
$ws.Range("A3").Value = "Tv do quarto"
$ws.Range("B3").NumberFormat = "#,##0"
$ws.Range("B3").HorizontalAlignment = -4152
$ws.Range("B3").Font.ThemeColor = 1
$ws.Range("C3").Value = 10
$ws.Range("C3").NumberFormat = "#,##0"
$ws.Range("C3").HorizontalAlignment = -4152
$ws.Range("C3").Font.Color = 0
$ws.Range("D3").Value = "'true"

# ---- Row 4: Tv da sala ----
$ws.Range("A4").Value = "Tv da sala"
$ws.Range("B4").Value = "Canal 3"
$ws.Range("B4").NumberFormat = "#,##0"
$ws.Range("B4").HorizontalAlignment = -4131
$ws.Range("B4").Font.ThemeColor = 1
$ws.Range("C4").Value = 10
$ws.Range("C4").NumberFormat = "#,##0"
$ws.Range("C4").HorizontalAlignment = -4152
$ws.Range("C4").Font.Color = 0
$ws.Range("D4").Value = "'true"

# ---- Row 5: Tv da cozinha ----
$ws.Range("A5").Value = "Tv da cozinha"
$ws.Range("B5").Value = "Canal 3"
$ws.Range("B5").NumberFormat = "#,##0"
$ws.Range("C5").Value = 10
$ws.Range("C5").NumberFormat = "#,##0"
$ws.Range("C5").HorizontalAlignment = -4152
$ws.Range("C5").Font.ThemeColor = 1
$ws.Range("D5").Value = "'true"
